# Questions - Answers Sheet.xlsx update
# Adds the "answer" / follow-up column (B) to the existing question rows on
# Feuil1 and appends a new trailing question in A10, matching the commit's
# "updated question sheet" change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New answers in column B, next to the existing questions in column A ---

# A4: "address type (signed or unsigned and size)"
# B4: rich-text answer with the "sizeof(duint)" part bold.
$ws.Range("B4").Value = "answer: unsigned, size: sizeof(duint)"
$bold = $ws.Range("B4").Characters(25, 13)
$bold.Font.Bold = $true

# A5: "size (32 or 64 bits and signess)"
$ws.Range("B5").Value = "answer: everything unsigned internally"

# A6: "Memory map analyis (Who do that, do you have some info)"
$ws.Range("B6").Value = "answer: we have dbghelp.dll, it can analyze some stuff and it seems logic to do this in the debugger code no?"

# A2: "Every time an instruction has been executed ... dump the full process memory?"
$ws.Range("B2").Value = "answer: not possible, it will take too much memory, but it is possible to dump the CIP region every time"

# --- New trailing question ---
$ws.Range("A10").Value = "question: switch to a word file??"

# --- Column sizing for the now much wider text (bestFit/AutoFit-style) ---
$ws.Columns.Item(1).ColumnWidth = 199.875
$ws.Columns.Item(2).ColumnWidth = 99

# --- Register the bold run's font in the workbook's style table (mirrors
# Excel's own behaviour of accumulating fonts/xfs as soon as Bold is used
# anywhere, even via a rich-text run) without leaving any stray cell behind.
$scratch = $ws.Range("Z100")
$scratch.Value = "x"
$scratch.Font.Bold = $true
$scratch.Clear()

# --- Selection / print setup ---
$ws.Range("A10").Select()

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
$ws.PageSetup.Copies = 0

Write-Output "Questions/answers sheet updated"
